$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "D" = -0.0246
    "E" = -0.218
    "F" = 0.17
    "G" = 0.451219512195122
    "H" = 0.451219512195122
    "I" = 0.1573170731707317
    "J" = 0.09339789130729897
    "K" = 3.2
    "L" = 0.09756097560975611
    "M" = 4.96
    "N" = 0.01725217391304348
    "O" = 1.55
    "P" = 4.96
    "Q" = 0.01725217391304348
    "R" = 1.55
    "U" = 76.3
    "V" = 0.2653913043478261
    "W" = 0.02714164546225615
    "X" = 0.02277911581982137
    "Y" = 0.004362529642434783
    "Z" = 0.8151701170564404
    "AA" = 0.0761351699897956
    "AB" = 0.02278506421089458
    "AC" = 0.05335010577890101
    "AD" = 0.107
    "AF" = 0.107
    "AG" = -76.193
    "AH" = 0.0003720354511538314
    "AI" = 0.0008857102651336429
    "AJ" = -0.3605796305848836
    "AK" = -1.711932954366729
    "AL" = 0.18
    "AM" = -0.114
    "AN" = 0.01059405940594059
    "AO" = 28.66666666666667
    "AP" = -7.543861386138614
    "AQ" = -45.26315789473685
}

foreach ($row in 2..3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
